# Writer now produces correct output for every FieldAttrType except TIME.
# Update the "Tasks" worksheet so rows DATE/TIMESTAMP/INTERVAL show green
# (done), row TIME keeps its yellow "not quite right" note updated for
# QlikView, and replace the old scratch notes with a new "Things to Do"
# list in column C (rows 13-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$GREEN = 5296274   # RGB(146, 208, 80) == FF92D050
$YELLOW = 65535     # RGB(255, 255, 0)  == FFFFFF00

# --- DATE row (8): writer now works, clear the old notes and make it green
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C8").Interior.Color = $GREEN

# --- TIME row (9): still the one type that doesn't work; keep yellow,
# clear the old "Reg-exp function" note, update the other note for QlikView
$ws.Range("B9").ClearContents()
$ws.Range("C9").Value = "Unsure of how to use Excel + QlikView to generate this data type"

# --- TIMESTAMP row (10): writer now works, clear notes, make green
$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C10").Interior.Color = $GREEN

# --- INTERVAL row (11): writer now works, clear notes (keep wrap on C11, just flip its color to green)
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("C11").Interior.Color = $GREEN

# --- Replace the single old scratch note (old B14) with a "Things to Do" list
$ws.Range("B14").Clear()

$ws.Range("C14").Value = "Timestamp format (figure out how to generate in QlikView)"
$ws.Range("C15").Value = "Reader + Writer: Null values should be allowed for QvxSpecialExtent, LocalDateTime, LocalDate, etc."
$ws.Range("C16").Value = "Deprecated verison of DateAndTimeCell"
$ws.Range("C17").Value = "Writer: Support for LocalDateTime, LocalDate, etc cell types"
$ws.Range("C18").Value = "Any dates before Feb 28, 1900 are one day off"

$ws.Range("C13").Value = "Things to Do"
$ws.Range("C13").Font.Bold = $true

$ws.Range("C15").Select()
